# Regenerate orders with updated distance/size codes.
# Mapping: D64 -> D69, D51 -> D55, D80 -> D86, S30 -> S31
# Applies to every text cell on the sheet (Condition, Filename_Left,
# Filename_Right, Distance, Size columns all contain these tokens).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$rowCount = $usedRange.Rows.Count
$colCount = $usedRange.Columns.Count

for ($r = 1; $r -le $rowCount; $r++) {
    for ($c = 1; $c -le $colCount; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value2

        if ($val -is [string]) {
            $newVal = $val.Replace("D64", "D69").Replace("D51", "D55").Replace("D80", "D86").Replace("S30", "S31")
            if ($newVal -ne $val) {
                $cell.Value2 = $newVal
            }
        }
    }
}
